# Updates cryptos list: refresh Price/Volume(1h) figures and swap the
# Kaspa/Polygon row order (rows 25-26), matching the GitHub Actions commit
# "Updated cryptos list on Fri Aug  2 17:34:33 UTC 2024 with GitHub Actions".
#
# Note: several Price values look like plain numbers (e.g. "559.59",
# "0.999") even though the sheet stores them as text. Setting .Value to such
# a string directly would make Excel coerce the cell to a Number, so for
# those we prefix with a leading apostrophe (forces text entry) and then
# reset NumberFormat/Style back to "Normal" so no visible formatting change
# is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.169.15"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "3.025.31"
$ws.Range("E3").Value = "  -2.49%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'559.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").Value = "'155.25"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").Value = "'0.564"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.12%  "
$ws.Range("D9").Value = "3.029.61"
$ws.Range("E9").Value = "  -2.14%  "
$ws.Range("D10").Value = "'0.114"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("D11").Value = "'6.44"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.91%  "
$ws.Range("D12").Value = "'0.368"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.13%  "
$ws.Range("D13").Value = "3.551.03"
$ws.Range("E13").Value = "  -2.35%  "
$ws.Range("E14").Value = "  -2.74%  "
$ws.Range("D15").Value = "63.211.35"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("D16").Value = "'24.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.47%  "
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "3.030.42"
$ws.Range("E18").Value = "  -2.03%  "
$ws.Range("D19").Value = "'398.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("D21").Value = "'12.04"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.17%  "
$ws.Range("E22").Value = "  -4.65%  "
$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("D24").Value = "'65.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.98%  "
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").Value = "'0.468"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.70%  "
$ws.Range("B26").Value = "Kaspa"
$ws.Range("C26").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D26").Value = "'0.190"
$ws.Range("D26").Style = "Normal"
$ws.Range("E27").Value = "  -1.50%  "
$ws.Range("D28").Value = "'8.71"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.81%  "
$ws.Range("E29").Value = "  -0.37%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  -0.43%  "
$ws.Range("D32").Value = "'20.48"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.56%  "
$ws.Range("D33").Value = "'159.96"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.58%  "
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("E35").Value = "  +2.10%  "
$ws.Range("E36").Value = "  -1.82%  "
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("D38").Value = "2.544.78"
$ws.Range("E38").Value = "  -5.79%  "
$ws.Range("E39").Value = "  -3.46%  "
$ws.Range("D40").Value = "'22.98"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.71%  "
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("D42").Value = "'37.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.34%  "
$ws.Range("D43").Value = "'0.671"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.28%  "
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("E45").Value = "  -0.97%  "
$ws.Range("E46").Value = "  -2.33%  "
$ws.Range("D47").Value = "'0.997"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.23%  "
$ws.Range("D48").Value = "'20.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.40%  "
$ws.Range("D49").Value = "'270.12"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.08%  "
$ws.Range("D50").Value = "'0.0948"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.31%  "
$ws.Range("D51").Value = "'10.47"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.03%  "
